$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text cells keep their exact literal representation (avoid Excel
# auto-converting numeric-looking strings like "589.08" into real numbers,
# which would drop significant trailing zeros / change representation).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.255.74"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -4.95%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.253.01"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "589.08"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -5.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.50"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -12.71%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.240.66"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -8.14%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.543"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -11.61%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -13.71%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.35%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -13.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.44"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -17.25%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -12.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.766.22"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.240.42"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -5.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "547.70"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -10.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.255.49"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -7.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.27"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -13.56%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -6.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.21"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -14.24%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -13.34%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -14.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.70"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -13.01%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -13.23%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -14.31%  "
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "29.48"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -12.67%  "
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.06"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -11.23%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -16.91%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -11.01%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -12.31%  "
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "548.99"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -14.00%  "
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.66"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -17.98%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -15.54%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0449"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -4.98%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "53.78"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -5.64%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0858"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -14.16%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -14.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.127"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -11.87%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.934.57"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -12.42%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -23.30%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -15.74%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₃0582"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -19.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "26.57"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -16.76%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -15.48%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "127.44"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -4.75%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -20.74%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -12.37%  "
